# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型"
# sheets to the refreshed values captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 649
$ws1.Range("F4").Value  = 223
$ws1.Range("F6").Value  = 9871
$ws1.Range("F7").Value  = 893
$ws1.Range("F9").Value  = 1238
$ws1.Range("F10").Value = 3934
$ws1.Range("F11").Value = 175
$ws1.Range("F12").Value = 115
$ws1.Range("F13").Value = 48
$ws1.Range("F15").Value = 286
$ws1.Range("F16").Value = 554
$ws1.Range("F18").Value = 269
$ws1.Range("F19").Value = 1464

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 649
$ws4.Range("F5").Value  = 223
$ws4.Range("F7").Value  = 9871
$ws4.Range("F8").Value  = 893
$ws4.Range("F10").Value = 1238
$ws4.Range("F11").Value = 3934
$ws4.Range("F12").Value = 175
$ws4.Range("F13").Value = 115
$ws4.Range("F14").Value = 48
$ws4.Range("F16").Value = 286
$ws4.Range("F17").Value = 554
$ws4.Range("F19").Value = 269
$ws4.Range("F20").Value = 1464
